$d = $word.ActiveDocument

$d.Content.Find.Execute("86÷7=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "92÷9=10, 2", 2) | Out-Null
$d.Content.Find.Execute("83÷6=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "66÷8=8, 2", 2) | Out-Null
$d.Content.Find.Execute("71÷5=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=14, 2", 2) | Out-Null
$d.Content.Find.Execute("32÷5=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=1, 4", 2) | Out-Null
$d.Content.Find.Execute("56÷4=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
$d.Content.Find.Execute("61÷9=6, 7", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=18, 1", 2) | Out-Null
$d.Content.Find.Execute("94÷8=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "27÷4=6, 3", 2) | Out-Null
$d.Content.Find.Execute("63÷9=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷7=12, 6", 2) | Out-Null
$d.Content.Find.Execute("80÷5=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=12, 0", 2) | Out-Null
$d.Content.Find.Execute("17÷3=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "87÷3=29, 0", 2) | Out-Null
$d.Content.Find.Execute("25÷7=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "80÷3=26, 2", 2) | Out-Null
$d.Content.Find.Execute("34÷5=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "37÷4=9, 1", 2) | Out-Null
$d.Content.Find.Execute("53÷7=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "71÷6=11, 5", 2) | Out-Null
$d.Content.Find.Execute("93÷4=23, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=13, 1", 2) | Out-Null
$d.Content.Find.Execute("76÷6=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=7, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷2=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=2, 2", 2) | Out-Null
$d.Content.Find.Execute("70÷9=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "73÷5=14, 3", 2) | Out-Null
$d.Content.Find.Execute("26÷7=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "23÷5=4, 3", 2) | Out-Null
$d.Content.Find.Execute("88÷9=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "18÷2=9, 0", 2) | Out-Null
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=5, 0", 2) | Out-Null
$d.Content.Find.Execute("91÷3=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "53÷9=5, 8", 2) | Out-Null
$d.Content.Find.Execute("40÷6=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "14÷6=2, 2", 2) | Out-Null
$d.Content.Find.Execute("19÷2=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷6=14, 4", 2) | Out-Null
$d.Content.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "96÷8=12, 0", 2) | Out-Null
$d.Content.Find.Execute("57÷4=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=12, 1", 2) | Out-Null
